$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking and must be forced to remain text
# (matching the original inline-string / text cell representation).
$textForceCells = @(
    'D4',
    'D5',
    'D6',
    'D7',
    'D9',
    'D10',
    'D14',
    'D15',
    'D16',
    'D18',
    'D20',
    'D21',
    'D23',
    'D24',
    'D26',
    'D27',
    'D33',
    'D34',
    'D35',
    'D36',
    'D38',
    'D43',
    'D44',
    'D45',
    'D47',
    'D50',
    'D51'
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values row by row.
$ws.Range('D2').Value = '51.063.27'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').Value = '2.951.85'
$ws.Range('E3').Value = '  -1.90%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '375.98'
$ws.Range('E5').Value = '  -2.58%  '
$ws.Range('D6').Value = '101.50'
$ws.Range('E6').Value = '  -3.75%  '
$ws.Range('D7').Value = '0.541'
$ws.Range('E7').Value = '  -1.49%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.593'
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('D10').Value = '36.46'
$ws.Range('E10').Value = '  -3.06%  '
$ws.Range('E11').Value = '  -0.74%  '
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').Value = '3.406.12'
$ws.Range('E13').Value = '  -2.15%  '
$ws.Range('D14').Value = '18.22'
$ws.Range('E14').Value = '  -1.85%  '
$ws.Range('D15').Value = '7.64'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').Value = '11.28'
$ws.Range('E16').Value = '  +50.36%  '
$ws.Range('D17').Value = '2.945.68'
$ws.Range('E17').Value = '  -1.84%  '
$ws.Range('D18').Value = '1.01'
$ws.Range('E18').Value = '  -1.81%  '
$ws.Range('D19').Value = '51.050.78'
$ws.Range('E19').Value = '  -1.53%  '
$ws.Range('D20').Value = '3.10'
$ws.Range('E20').Value = '  -6.55%  '
$ws.Range('D21').Value = '12.53'
$ws.Range('E21').Value = '  -4.04%  '
$ws.Range('D22').Value = '0.0₃0958'
$ws.Range('E22').Value = '  -1.20%  '
$ws.Range('D23').Value = '266.37'
$ws.Range('D24').Value = '68.88'
$ws.Range('E24').Value = '  -0.69%  '
$ws.Range('E25').Value = '  +6.53%  '
$ws.Range('D26').Value = '8.19'
$ws.Range('E26').Value = '  -2.76%  '
$ws.Range('D27').Value = '7.65'
$ws.Range('E27').Value = '  -0.25%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('E30').Value = '  -4.53%  '
$ws.Range('E31').Value = '  -4.14%  '
$ws.Range('E32').Value = '  +0.66%  '
$ws.Range('D33').Value = '50.88'
$ws.Range('E33').Value = '  -0.67%  '
$ws.Range('D34').Value = '2.06'
$ws.Range('E34').Value = '  -0.69%  '
$ws.Range('D35').Value = '33.59'
$ws.Range('E35').Value = '  -4.95%  '
$ws.Range('D36').Value = '0.0445'
$ws.Range('E36').Value = '  -2.67%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').Value = '3.19'
$ws.Range('E38').Value = '  +3.60%  '
$ws.Range('E39').Value = '  -1.06%  '
$ws.Range('E40').Value = '  -4.08%  '
$ws.Range('E41').Value = '  -3.00%  '
$ws.Range('E42').Value = '  -4.90%  '
$ws.Range('D43').Value = '120.79'
$ws.Range('E43').Value = '  -1.69%  '
$ws.Range('D44').Value = '21.42'
$ws.Range('E44').Value = '  -2.46%  '
$ws.Range('D45').Value = '3.42'
$ws.Range('E45').Value = '  +2.59%  '
$ws.Range('E46').Value = '  -1.16%  '
$ws.Range('D47').Value = '0.272'
$ws.Range('E47').Value = '  -2.80%  '
$ws.Range('E48').Value = '  -2.00%  '
$ws.Range('D49').Value = '1.995.82'
$ws.Range('E49').Value = '  -2.54%  '
$ws.Range('D50').Value = '0.0329'
$ws.Range('E50').Value = '  -2.32%  '
$ws.Range('D51').Value = '1.31'
$ws.Range('E51').Value = '  +1.88%  '

# Remove the temporary text-number-format so cell styles match the original
# (unstyled) inline-string cells exactly.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).ClearFormats()
}
